# Apply "diagram labeling thru 2016" edits:
#  1. Update a handful of H-column (subcategory) labels.
#  2. Remove the now-unused I column ("is_viewed") entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Relabel subcategory values in column H.
$updates = @{
    11 = "photo(s)"
    12 = "photo(s)"
    13 = "drawing(s)"
    17 = "bar chart(s)"
    18 = "bar chart(s)"
    35 = "line graph(s)"
    39 = "data display"
    45 = "bar chart(s)"
    46 = "line graph(s)"
    47 = "bar chart(s)"
    48 = "line graph(s)"
    56 = "data collection, data analysis, data gathering diagram"
    60 = "data collection, data analysis, data gathering diagram"
    62 = "data collection, data analysis, data gathering diagram"
    64 = "line graph(s)"
    67 = "bar chart(s)"
    70 = "line graph(s)"
}

foreach ($row in $updates.Keys) {
    $ws.Range("H$row").Value = $updates[$row]
}

# 2. Delete column I ("is_viewed") entirely, shrinking the used range to A1:H70.
$ws.Columns("I:I").Delete()
